$wb = $excel.ActiveWorkbook

# Rename the sheets
$wsZoo = $wb.Worksheets.Item(1)
$wsExpert = $wb.Worksheets.Item(2)
$wsZoo.Name = "ZooniverseTags"
$wsExpert.Name = "ExpertTags"

# Insert a new first row on the "ExpertTags" sheet with the tag "name"
# (shifting all existing rows down by one)
$wsExpert.Rows.Item(1).Insert()
$wsExpert.Range("A1").Value = "name"
$wsExpert.Rows.Item(1).RowHeight = 15.75
